# added 4wk low sales check
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2 (W10)
$wsForecast.Range("D2").Value = 5
$wsForecast.Range("H2").Value = 0.17
$wsForecast.Range("L2").Value = 1.09

# Row 3 (W11)
$wsForecast.Range("L3").Value = 0.9399999999999999

# Row 4 (W12)
$wsForecast.Range("L4").Value = 0.93

# Row 5 (W13)
$wsForecast.Range("L5").Value = 1.08

# Row 6 (W14)
$wsForecast.Range("L6").Value = 1.18

# Row 7 (W15)
$wsForecast.Range("L7").Value = 1.15

# Row 8 (W16)
$wsForecast.Range("L8").Value = 1.09

# Row 9 (W17)
$wsForecast.Range("L9").Value = 1.01

# Row 10 (W18)
$wsForecast.Range("D10").Value = 4
$wsForecast.Range("L10").Value = 0.91

# Row 11 (W19)
$wsForecast.Range("D11").Value = 4
$wsForecast.Range("L11").Value = 0.91

# Row 12 (W20)
$wsForecast.Range("D12").Value = 4
$wsForecast.Range("L12").Value = 0.91

# Row 13 (W21)
$wsForecast.Range("D13").Value = 4
$wsForecast.Range("L13").Value = 1.1

# Row 14 (W22)
$wsForecast.Range("D14").Value = 4
$wsForecast.Range("L14").Value = 1.19

# Row 15 (W23)
$wsForecast.Range("L15").Value = 1.2

# Row 16 (W24)
$wsForecast.Range("L16").Value = 1.01

# Row 17 (W25)
$wsForecast.Range("L17").Value = 1.08

# --- Summary sheet updates ---
# These "Value" cells hold numeric-looking text (matching the source
# report's inline-string formatting), so force text entry the same way a
# user would in Excel -- with a leading apostrophe -- rather than letting
# the numeric strings be auto-converted to Number cells.

$wsSummary.Range("B9").Value = "'79"
$wsSummary.Range("B10").Value = "'44"
$wsSummary.Range("B11").Value = "'23"
$wsSummary.Range("B14").Value = "'4"
